# "update heur. for v2; start v3"
# Rename the v2 heuristic-named sheets, drop the scratch Sheet6, and move
# the live selection on to Randomish2 (v3) while updating each sheet's
# remembered cursor position.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Rename sheets for v2 -------------------------------------------------
$wb.Worksheets("StronglySkewed").Name = "Skewed1"
$wb.Worksheets("FrontLoaded").Name = "Skewed2"
$wb.Worksheets("Randomish").Name = "Randomish1"
# "Randomish2" keeps its name.

# --- Drop the empty scratch sheet ------------------------------------------
$wb.Worksheets("Sheet6").Delete()

# --- Skewed1 (was StronglySkewed): scrolled down, selection unchanged -----
$ws = $wb.Worksheets("Skewed1")
$ws.Activate()
$ws.Range("E61").Select()

# --- Skewed2 (was FrontLoaded): scrolled down, new selection --------------
$ws = $wb.Worksheets("Skewed2")
$ws.Activate()
$ws.Range("F74").Select()

# --- Randomish1 (was Randomish): new selection -----------------------------
$ws = $wb.Worksheets("Randomish1")
$ws.Activate()
$ws.Range("O36").Select()

# --- Randomish2: new selection, becomes the active/selected tab (start v3) -
$ws = $wb.Worksheets("Randomish2")
$ws.Activate()
$ws.Range("P44").Select()
